$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.973936579530329
$ws.Cells.Item(2, 4).Value = 2.888251324519063
$ws.Cells.Item(2, 5).Value = 16.71456242210408
$ws.Cells.Item(2, 6).Value = 22.91516405752209
$ws.Cells.Item(2, 7).Value = 3.554465621914874
$ws.Cells.Item(2, 9).Value = 17.51723746316919
$ws.Cells.Item(2, 14).Value = 17.13588176064243
$ws.Cells.Item(2, 15).Value = 19.41803292077183
$ws.Cells.Item(3, 2).Value = 7.726105456002573
$ws.Cells.Item(3, 4).Value = 2.868465467081656
$ws.Cells.Item(3, 5).Value = 15.74734985263152
$ws.Cells.Item(3, 6).Value = 22.18566614097691
$ws.Cells.Item(3, 7).Value = 3.558084161067144
$ws.Cells.Item(3, 9).Value = 17.43545559361282
$ws.Cells.Item(3, 14).Value = 16.94880368210766
$ws.Cells.Item(3, 15).Value = 18.92251440240058
$ws.Cells.Item(4, 2).Value = 7.569968192699863
$ws.Cells.Item(4, 4).Value = 2.856048993626553
$ws.Cells.Item(4, 5).Value = 15.12794543687676
$ws.Cells.Item(4, 6).Value = 21.73354054194773
$ws.Cells.Item(4, 7).Value = 3.560419138670963
$ws.Cells.Item(4, 9).Value = 17.38881070702378
$ws.Cells.Item(4, 14).Value = 16.8348163882461
$ws.Cells.Item(4, 15).Value = 18.61841796641237
$ws.Cells.Item(5, 2).Value = 7.505429865115993
$ws.Cells.Item(5, 4).Value = 2.850921696415942
$ws.Cells.Item(5, 5).Value = 14.86937712186838
$ws.Cells.Item(5, 6).Value = 21.54856223081956
$ws.Cells.Item(5, 7).Value = 3.561399240755241
$ws.Cells.Item(5, 9).Value = 17.37072111760673
$ws.Cells.Item(5, 14).Value = 16.78863424108822
$ws.Cells.Item(5, 15).Value = 18.49473964450635
$ws.Cells.Item(6, 2).Value = 7.494660968687292
$ws.Cells.Item(6, 4).Value = 2.85006623315215
$ws.Cells.Item(6, 5).Value = 14.8260790259931
$ws.Cells.Item(6, 6).Value = 21.51781232752263
$ws.Cells.Item(6, 7).Value = 3.561563715616634
$ws.Cells.Item(6, 9).Value = 17.36777342702478
$ws.Cells.Item(6, 14).Value = 16.78098333545059
$ws.Cells.Item(6, 15).Value = 18.47422376162881
$ws.Cells.Item(7, 2).Value = 7.569101375577702
$ws.Cells.Item(7, 4).Value = 2.855980118445007
$ws.Cells.Item(7, 5).Value = 15.12448283059672
$ws.Cells.Item(7, 6).Value = 21.73104838485154
$ws.Cells.Item(7, 7).Value = 3.560432240795761
$ws.Cells.Item(7, 9).Value = 17.38856299825989
$ws.Cells.Item(7, 14).Value = 16.8341924095864
$ws.Cells.Item(7, 15).Value = 18.61674873665768
$ws.Cells.Item(8, 2).Value = 7.889357710724819
$ws.Cells.Item(8, 4).Value = 2.881484681021925
$ws.Cells.Item(8, 5).Value = 16.3865093745028
$ws.Cells.Item(8, 6).Value = 22.66469596995105
$ws.Cells.Item(8, 7).Value = 3.555689875988606
$ws.Cells.Item(8, 9).Value = 17.48831115075855
$ws.Cells.Item(8, 14).Value = 17.07122325684954
$ws.Cells.Item(8, 15).Value = 19.24726691713163
$ws.Cells.Item(9, 2).Value = 8.4824544463425
$ws.Cells.Item(9, 4).Value = 2.929380739128483
$ws.Cells.Item(9, 5).Value = 18.80736662390711
$ws.Cells.Item(9, 6).Value = 24.4489572850658
$ws.Cells.Item(9, 7).Value = 3.547282679132513
$ws.Cells.Item(9, 9).Value = 17.71126666152012
$ws.Cells.Item(9, 14).Value = 17.54092831277462
$ws.Cells.Item(9, 15).Value = 20.47631796316274
$ws.Cells.Item(10, 2).Value = 8.892805469389883
$ws.Cells.Item(10, 4).Value = 2.96329102212657
$ws.Cells.Item(10, 5).Value = 20.49547908228563
$ws.Cells.Item(10, 6).Value = 25.71537767688056
$ws.Cells.Item(10, 7).Value = 3.541642386682013
$ws.Cells.Item(10, 9).Value = 17.89037526883937
$ws.Cells.Item(10, 14).Value = 17.88621396499388
$ws.Cells.Item(10, 15).Value = 21.36392112697861
$ws.Cells.Item(11, 2).Value = 9.073239364152236
$ws.Cells.Item(11, 4).Value = 2.978442842350302
$ws.Cells.Item(11, 5).Value = 21.22102878584207
$ws.Cells.Item(11, 6).Value = 26.2789289706539
$ws.Cells.Item(11, 7).Value = 3.539191319903617
$ws.Cells.Item(11, 9).Value = 17.97487022382017
$ws.Cells.Item(11, 14).Value = 18.04275315944065
$ws.Cells.Item(11, 15).Value = 21.76229304255146
$ws.Cells.Item(12, 2).Value = 9.140615808878705
$ws.Cells.Item(12, 4).Value = 2.984141199598102
$ws.Cells.Item(12, 5).Value = 21.48972210935719
$ws.Cells.Item(12, 6).Value = 26.4903152199852
$ws.Cells.Item(12, 7).Value = 3.538279535533335
$ws.Cells.Item(12, 9).Value = 18.0072731922733
$ws.Cells.Item(12, 14).Value = 18.10190618682539
$ws.Cells.Item(12, 15).Value = 21.91221204852734
$ws.Cells.Item(13, 2).Value = 9.126148118074333
$ws.Cells.Item(13, 4).Value = 2.982915702787114
$ws.Cells.Item(13, 5).Value = 21.43212296421097
$ws.Cells.Item(13, 6).Value = 26.44488234123333
$ws.Cells.Item(13, 7).Value = 3.538475177899426
$ws.Cells.Item(13, 9).Value = 18.00027699228527
$ws.Cells.Item(13, 14).Value = 18.08917283935257
$ws.Cells.Item(13, 15).Value = 21.87996832303193
$ws.Cells.Item(14, 2).Value = 9.078801786399685
$ws.Cells.Item(14, 4).Value = 2.978912436746874
$ws.Cells.Item(14, 5).Value = 21.24325551487777
$ws.Cells.Item(14, 6).Value = 26.29636144028487
$ws.Cells.Item(14, 7).Value = 3.539115979214815
$ws.Cells.Item(14, 9).Value = 17.97752804593542
$ws.Cells.Item(14, 14).Value = 18.04762249193304
$ws.Cells.Item(14, 15).Value = 21.77464658914564
$ws.Cells.Item(15, 2).Value = 9.049675604494444
$ws.Cells.Item(15, 4).Value = 2.976455205106719
$ws.Cells.Item(15, 5).Value = 21.12678111820752
$ws.Cells.Item(15, 6).Value = 26.20511928788498
$ws.Cells.Item(15, 7).Value = 3.539510618406177
$ws.Cells.Item(15, 9).Value = 17.9636457851309
$ws.Cells.Item(15, 14).Value = 18.02215401479316
$ws.Cells.Item(15, 15).Value = 21.71000761895239
$ws.Cells.Item(16, 2).Value = 8.880883512103345
$ws.Cells.Item(16, 4).Value = 2.962295308117491
$ws.Cells.Item(16, 5).Value = 20.44721311684653
$ws.Cells.Item(16, 6).Value = 25.67827694917971
$ws.Cells.Item(16, 7).Value = 3.541804868354967
$ws.Cells.Item(16, 9).Value = 17.8849117193195
$ws.Cells.Item(16, 14).Value = 17.87596887941633
$ws.Cells.Item(16, 15).Value = 21.33776399663439
$ws.Cells.Item(17, 2).Value = 8.775697920554309
$ws.Cells.Item(17, 4).Value = 2.9535385342362
$ws.Cells.Item(17, 5).Value = 20.01949224721977
$ws.Cells.Item(17, 6).Value = 25.3517026014505
$ws.Cells.Item(17, 7).Value = 3.543241618569027
$ws.Cells.Item(17, 9).Value = 17.83736449730081
$ws.Cells.Item(17, 14).Value = 17.78611779882749
$ws.Cells.Item(17, 15).Value = 21.10790372692033
$ws.Cells.Item(18, 2).Value = 8.714613687309544
$ws.Cells.Item(18, 4).Value = 2.948476059198042
$ws.Cells.Item(18, 5).Value = 19.76948715586678
$ws.Cells.Item(18, 6).Value = 25.16269573567658
$ws.Cells.Item(18, 7).Value = 3.544078805245713
$ws.Cells.Item(18, 9).Value = 17.81030229408061
$ws.Cells.Item(18, 14).Value = 17.73439049486328
$ws.Cells.Item(18, 15).Value = 20.9751943103584
$ws.Cells.Item(19, 2).Value = 8.693832984641938
$ws.Cells.Item(19, 4).Value = 2.946757555924558
$ws.Cells.Item(19, 5).Value = 19.68415305554162
$ws.Cells.Item(19, 6).Value = 25.09850720800331
$ws.Cells.Item(19, 7).Value = 3.544364121757167
$ws.Cells.Item(19, 9).Value = 17.80118939117589
$ws.Cells.Item(19, 14).Value = 17.71686987832943
$ws.Cells.Item(19, 15).Value = 20.93018037350943
$ws.Cells.Item(20, 2).Value = 8.786955979029495
$ws.Cells.Item(20, 4).Value = 2.954473378242521
$ws.Cells.Item(20, 5).Value = 20.06543668049255
$ws.Cells.Item(20, 6).Value = 25.38658970316626
$ws.Cells.Item(20, 7).Value = 3.543087556526535
$ws.Cells.Item(20, 9).Value = 17.84239660774314
$ws.Cells.Item(20, 14).Value = 17.79568785271765
$ws.Cells.Item(20, 15).Value = 21.13242566986916
$ws.Cells.Item(21, 2).Value = 9.092734726993216
$ws.Cells.Item(21, 4).Value = 2.980089358902621
$ws.Cells.Item(21, 5).Value = 21.29889453732582
$ws.Cells.Item(21, 6).Value = 26.34004199921628
$ws.Cells.Item(21, 7).Value = 3.538927316677012
$ws.Cells.Item(21, 9).Value = 17.9841991494933
$ws.Cells.Item(21, 14).Value = 18.05983060728496
$ws.Cells.Item(21, 15).Value = 21.8056087360082
$ws.Cells.Item(22, 2).Value = 9.287021218569393
$ws.Cells.Item(22, 4).Value = 2.996601823451826
$ws.Cells.Item(22, 5).Value = 22.06974545028082
$ws.Cells.Item(22, 6).Value = 26.95133001253828
$ws.Cells.Item(22, 7).Value = 3.536303788819382
$ws.Cells.Item(22, 9).Value = 18.07923203987225
$ws.Cells.Item(22, 14).Value = 18.23171317604655
$ws.Cells.Item(22, 15).Value = 22.24005346421713
$ws.Cells.Item(23, 2).Value = 9.183851410766957
$ws.Cells.Item(23, 4).Value = 2.987809729864144
$ws.Cells.Item(23, 5).Value = 21.66154396560651
$ws.Cells.Item(23, 6).Value = 26.62622303281883
$ws.Cells.Item(23, 7).Value = 3.537695321746732
$ws.Cells.Item(23, 9).Value = 18.02830472145926
$ws.Cells.Item(23, 14).Value = 18.14006050298211
$ws.Cells.Item(23, 15).Value = 22.00873637834044
$ws.Cells.Item(24, 2).Value = 8.781868111808409
$ws.Cells.Item(24, 4).Value = 2.954050822766922
$ws.Cells.Item(24, 5).Value = 20.0446779839577
$ws.Cells.Item(24, 6).Value = 25.37082113682247
$ws.Cells.Item(24, 7).Value = 3.543157173144086
$ws.Cells.Item(24, 9).Value = 17.84012073718833
$ws.Cells.Item(24, 14).Value = 17.79136144795663
$ws.Cells.Item(24, 15).Value = 21.12134103057043
$ws.Cells.Item(25, 2).Value = 8.326196084739937
$ws.Cells.Item(25, 4).Value = 2.916649087185314
$ws.Cells.Item(25, 5).Value = 18.14787061189909
$ws.Cells.Item(25, 6).Value = 23.97300802016276
$ws.Cells.Item(25, 7).Value = 3.54946228545465
$ws.Cells.Item(25, 9).Value = 17.64816122966617
$ws.Cells.Item(25, 14).Value = 17.4136265380636
$ws.Cells.Item(25, 15).Value = 20.14578213644698
